$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "Android Compose betanulás "
$ws.Range("C12").Value = "Android kliens -  kosár, checkout képernyők megvalósítása,  termék filterezés megvalósítása - Diplomamunka írás"
$ws.Range("C11").Value = "Android kliens - termékek, termékkategóriák, autentikáció megvalósítás, lokális adatbázis megvalósítás"
$ws.Range("C10").Value = "Android kliens - architektura, app skeleton kialakítása, technológiák kiválasztása, függőségek hozzáadása, navigáció megvalósítása"
$ws.Range("C9").Value = "Android Architecture, Navigation betanulás"

$ws.Range("C13").Select()
